# Added sections & flags. Omitted type conversions in .tt file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# A1 "register" and B1 "data" keep their text; D1 renamed group -> section
$ws.Cells.Item(1, 1).Value = "register"
$ws.Cells.Item(1, 2).Value = "data"
$ws.Cells.Item(1, 3).Value = "isFlag"
$ws.Cells.Item(1, 4).Value = "section"

# --- Data rows ----------------------------------------------------------
# Column A: register/flag identifiers (string), Column B: register data /
# ON-OFF flag state. Columns C (formula) and D (section number) unchanged.

$ws.Cells.Item(2, 1).Value = "F[20:Flg_Test]"
$ws.Cells.Item(2, 2).Value = "ON"

$ws.Cells.Item(3, 1).Value = "F[21:Flg_Test]"
$ws.Cells.Item(3, 2).Value = "OFF"

$ws.Cells.Item(4, 1).Value = "F[22:Flg_Test]"
$ws.Cells.Item(4, 2).Value = "OFF"

$ws.Cells.Item(5, 1).Value = "F[23:FLg_Test]"
$ws.Cells.Item(5, 2).Value = "ON"

$ws.Cells.Item(6, 1).Value = "R[1:Reg_Test]"
$ws.Cells.Item(6, 2).Value = 30.5

$ws.Cells.Item(7, 1).Value = "R[2:Reg_Test]"
$ws.Cells.Item(7, 2).Value = 300

$ws.Cells.Item(8, 1).Value = "R[3:Reg_Test]"
$ws.Cells.Item(8, 2).Value = 30.7

$ws.Cells.Item(9, 1).Value = "R[4:Reg_Test]"
$ws.Cells.Item(9, 2).Value = 30.8

$ws.Cells.Item(10, 1).Value = "R[5:Reg_Test]"
$ws.Cells.Item(10, 2).Value = 30.9

$ws.Cells.Item(11, 1).Value = "R[6:Reg_Test]"
$ws.Cells.Item(11, 2).Value = 31

$ws.Cells.Item(12, 1).Value = "R[7:Reg_Test]"
$ws.Cells.Item(12, 2).Value = 31.1

$ws.Cells.Item(13, 1).Value = "R[8:Reg_Test]"
$ws.Cells.Item(13, 2).Value = 31.2

$ws.Cells.Item(14, 1).Value = "R[9:Reg_Test]"
$ws.Cells.Item(14, 2).Value = 31.3

$ws.Cells.Item(15, 1).Value = "R[10:Reg_Test]"
$ws.Cells.Item(15, 2).Value = 31.4

$ws.Cells.Item(16, 1).Value = "R[11:Reg_Test]"
$ws.Cells.Item(16, 2).Value = 31.5

$ws.Cells.Item(17, 1).Value = "R[12:Reg_Test]"
$ws.Cells.Item(17, 2).Value = 31.6

$ws.Cells.Item(18, 1).Value = "R[13:Reg_Test]"
$ws.Cells.Item(18, 2).Value = 31.7

$ws.Cells.Item(19, 1).Value = "R[14:Reg_Test]"
$ws.Cells.Item(19, 2).Value = 31.8

$ws.Cells.Item(20, 1).Value = "R[15:Reg_Test]"
$ws.Cells.Item(20, 2).Value = 31.9

$ws.Cells.Item(21, 1).Value = "R[16:Reg_Test]"
$ws.Cells.Item(21, 2).Value = 32

# --- Column widths --------------------------------------------------
# Target widths are 16.25 (col A) and 9.75 (cols B:C) characters. The
# host only persists ColumnWidth at whole-pixel granularity (7 px per
# character + 5 px padding, then rounded to the nearest pixel), so the
# literal values 16.25 / 9.75 would actually round away to 17 / 10.43.
# Feeding in the character widths whose resulting pixel width is
# closest to the intended 16.25 / 9.75 gets us the nearest achievable
# stored width (16.29 / 9.71).
$ws.Columns.Item(1).ColumnWidth = 15.57
$ws.Range("B1:C1").ColumnWidth = 9.0

# --- Selection --------------------------------------------------------
$ws.Range("H11").Select()
